$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look numeric need an explicit Text format
# so Excel stores them verbatim (e.g. "1.00", "3.62") instead of coercing
# them to a Double and losing the original formatting.
$textFormatCells = @(
    'D5',
    'D6',
    'D8',
    'D11',
    'D15',
    'D17',
    'D18',
    'D22',
    'D23',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D41',
    'D44',
    'D45',
    'D46',
    'D48',
    'D51'
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Updated coin prices / 1h volume percentages scraped this run,
# plus the Maker/Stacks row swap and the OceanProtocol -> Monero row.
$cellUpdates = [ordered]@{
    'D2' = '69.370.83'
    'E2' = '  -1.99%  '
    'D3' = '3.491.43'
    'E3' = '  -1.72%  '
    'E4' = '  -0.08%  '
    'D5' = '611.88'
    'E5' = '  +5.20%  '
    'D6' = '186.08'
    'E6' = '  +0.48%  '
    'E7' = '  +0.39%  '
    'D8' = '1.00'
    'E8' = '  -0.06%  '
    'E9' = '  -1.74%  '
    'E10' = '  +0.31%  '
    'D11' = '53.14'
    'E11' = '  -2.24%  '
    'E12' = '  -1.79%  '
    'E13' = '  +0.81%  '
    'D14' = '4.041.70'
    'E14' = '  -1.85%  '
    'D15' = '606.75'
    'E15' = '  +6.18%  '
    'D16' = '69.377.80'
    'E16' = '  -1.98%  '
    'D17' = '12.62'
    'E17' = '  +2.36%  '
    'D18' = '18.85'
    'E18' = '  -1.90%  '
    'D19' = '3.489.12'
    'E19' = '  -0.55%  '
    'E20' = '  -0.09%  '
    'E21' = '  -1.36%  '
    'D22' = '17.17'
    'E22' = '  -2.93%  '
    'D23' = '105.90'
    'E23' = '  +10.70%  '
    'E24' = '  +2.37%  '
    'E25' = '  +1.61%  '
    'E26' = '  +2.73%  '
    'D27' = '10.95'
    'E27' = '  -2.48%  '
    'D28' = '9.83'
    'E28' = '  +7.43%  '
    'D29' = '33.64'
    'E29' = '  +3.59%  '
    'D30' = '6.96'
    'E30' = '  -3.29%  '
    'D31' = '12.36'
    'E31' = '  +0.89%  '
    'E32' = '  -1.05%  '
    'D33' = '3.91'
    'E33' = '  +14.92%  '
    'D34' = '63.13'
    'E34' = '  -0.24%  '
    'D35' = '3.20'
    'E35' = '  -6.27%  '
    'D36' = '0.999'
    'E36' = '  -0.13%  '
    'D37' = '519.91'
    'E37' = '  -4.58%  '
    'D38' = '0.396'
    'E38' = '  -3.90%  '
    'B39' = 'Stacks'
    'C39' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D39' = '3.62'
    'E39' = '  +5.86%  '
    'B40' = 'Maker'
    'C40' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D40' = '3.587.11'
    'E40' = '  +0.32%  '
    'D41' = '36.78'
    'E41' = '  -2.66%  '
    'E42' = '  -2.99%  '
    'E43' = '  +0.97%  '
    'D44' = '0.0462'
    'E44' = '  +1.94%  '
    'D45' = '2.95'
    'E45' = '  +1.18%  '
    'D46' = '0.142'
    'E46' = '  +3.37%  '
    'E47' = '  -4.02%  '
    'D48' = '8.83'
    'E48' = '  -4.85%  '
    'E49' = '  +0.42%  '
    'B51' = 'Monero'
    'C51' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D51' = '130.85'
    'E51' = '  -2.46%  '
}
foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}
